# Swap the data content of rows 2 and 3 on the active sheet.
# (The edit author re-ordered two observation records; every populated
# cell in row 2 and row 3 trades places, including cells that are only
# present in one of the two rows.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2 <-> A3
$ws.Range("A2").Value = 102077473
$ws.Range("A3").Value = 80139019

# B2 <-> B3
$ws.Range("B2").Value = 96367
$ws.Range("B3").Value = 88921

# E2 <-> E3
$ws.Range("E2").Value = 219874
$ws.Range("E3").Value = 5741

# F2 <-> F3
$ws.Range("F2").Formula = "'" + 'Nattviol'
$ws.Range("F3").Formula = "'" + 'Tjockfotad fingersvamp'

# G2 <-> G3
$ws.Range("G2").Formula = "'" + 'Platanthera bifolia'
$ws.Range("G3").Formula = "'" + 'Ramaria flavescens'

# H2 <-> H3
$ws.Range("H2").Formula = "'" + '(L.) Rich.'
$ws.Range("H3").Formula = "'" + '(Schaeff.) R. H. Petersen'

# I2 <-> I3
$ws.Range("I2").Formula = "'" + '3'
$ws.Range("I3").Formula = "'" + ''

# J2 <-> J3
$ws.Range("J2").Formula = "'" + ''
$ws.Range("J3").ClearContents()

# K2 <-> K3
$ws.Range("K2").Formula = "'" + ''
$ws.Range("K3").ClearContents()

# L2 <-> L3
$ws.Range("L2").Formula = "'" + ''
$ws.Range("L3").ClearContents()

# N2 <-> N3
$ws.Range("N2").Formula = "'" + ''
$ws.Range("N3").ClearContents()

# P2 <-> P3
$ws.Range("P2").Formula = "'" + 'Tvetaspåret, Tveta, Srm'
$ws.Range("P3").Formula = "'" + 'Tveta friluftsgård, 300 m V om, Srm'

# Q2 <-> Q3
$ws.Range("Q2").Value = 647720.9098417715
$ws.Range("Q3").Value = 648222.682956806

# R2 <-> R3
$ws.Range("R2").Value = 6560694.968483768
$ws.Range("R3").Value = 6560420.292955686

# S2 <-> S3
$ws.Range("S2").Value = 10
$ws.Range("S3").Value = 50

# Y2 <-> Y3
$ws.Range("Y2").Formula = "'" + '2022-06-28'
$ws.Range("Y3").Formula = "'" + '2019-09-27'

# AA2 <-> AA3
$ws.Range("AA2").Formula = "'" + '2022-07-05'
$ws.Range("AA3").Formula = "'" + '2019-09-27'

# AF2 <-> AF3
$ws.Range("AF2").Formula = "'" + ''
$ws.Range("AF3").ClearContents()

# AI2 <-> AI3
$ws.Range("AI2").ClearContents()
$ws.Range("AI3").Formula = "'" + 'barrskog'

# AW2 <-> AW3
$ws.Range("AW2").Formula = "'" + 'Åsa Johansson'
$ws.Range("AW3").Formula = "'" + 'Hans Rydberg'

# AX2 <-> AX3
$ws.Range("AX2").Formula = "'" + 'Åsa Johansson'
$ws.Range("AX3").Formula = "'" + 'Hans Rydberg'

